# Update column G ("K" - strikeouts) values for rows 2-14 per regenerated save_data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 3
    4  = 0
    5  = 0
    6  = 2
    7  = 3
    8  = 5
    9  = 1
    10 = 3
    11 = 1
    12 = 0
    13 = 1
    14 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
